$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.555.36"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.873.77"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4730"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2919"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("E13").Value = "  +4.94%  "
$ws.Range("D14").Value = "1.869.51"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.160"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "30.522.84"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007522"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "2.114.24"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9990"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.264"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.175"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.920"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09992"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.350"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.510"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.293"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04809"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6962"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.744"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.231"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.967"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4191"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8342"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.372"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.982"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "927.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05645"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.35%  "
